# The deck's theme parts (ppt/theme/theme1.xml and ppt/theme/theme2.xml) had
# their "Office Theme" / "Integral" color schemes swapped: theme2.xml (the
# part actually wired to the slide master / used by every slide) needs to go
# from the "Red Violet" / Integral palette to the plain "Office" palette that
# used to live in theme1.xml.
#
# The PowerPoint object model doesn't give us a way to rewrite a theme part's
# raw XML or to re-point a master at a different theme file, but it does let
# us edit every slot of the *current* theme color scheme in place (the 12
# DrawingML theme colors: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
# Rewriting those 12 slots on the live theme reproduces the visible effect of
# the swap - every slide, layout and the slide master itself switch from the
# pink/purple "Integral" palette to the blue/orange "Office" palette.

function ConvertTo-RgbLong($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint's ColorFormat.RGB is a VBA-style Long: 0xBBGGRR
    return ($b * 65536) + ($g * 256) + $r
}

# Target values, taken from the "Office Theme" palette (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink), in ThemeColorScheme.Colors(1..12) order.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation

# Any slide exposes the presentation's single live ThemeColorScheme; slide 1
# is as good an anchor as any.
$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $themeColors.Colors($i).RGB = ConvertTo-RgbLong $officeThemeColors[$i - 1]
}
